$d = $word.ActiveDocument

# 1) Replace the title paragraph text.
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Text = "我的青春年華，宛如含苞待放的玫瑰，在時光流轉中，緩緩地綻放著。"

# 2) Replace the "invention" paragraph (simplified -> traditional, rewritten).
$p4 = $d.Paragraphs.Item(4)
$p4.Range.Text = "發明，是人類智慧的火花，照亮了文明的長河。它就像風箏般，牽引著人類進步的絲線，越飛越高，越飛越遠。發明家們，是那些將夢想化作現實的魔法師，他們用靈魂的火焰，點燃了物質世界的潛能。他們的靈感，如夜空中閃耀的星辰，指引著人類走向更美好的未來。然而，發明之路，也充滿了荊棘與坎坷。這些先驅者，常常被世俗的迷霧所遮蔽，他們的努力，往往得不到應有的回報。但當他們看到自己的發明，化作世間的福祉，照亮了無數人的生命時，他們的心中，便充滿了無比的欣慰與滿足。沒有發明家，人類便如困在黑暗中的迷途羔羊，無助地徘徊在生存的邊緣。是他們，用智慧與毅力，為我們創造了光明與希望，讓人類文明得以延續，生生不息。"

# 3) Replace the "I, laughable to say..." paragraph.
$p5 = $d.Paragraphs.Item(5)
$p5.Range.Text = "這一生，我已嚐遍了無與倫比的喜悅，多到我這些年的生命幾乎浸淫在無止境的狂喜之中。人們說我勤奮，也許吧，如果思考也等同於勞動，因為我幾乎將所有清醒時光都奉獻給了思想。但若說工作是按著嚴格的規矩，在特定的時間內完成特定的表現，那我恐怕是最不濟的懶人。任何強迫性的努力都需要犧牲生命力，而我從未願意付出這樣的代價。相反地，我的思想讓我如沐春風，在思想的海洋裡自由翱翔，那是生命最美妙的旋律，也是我生命的全部。"

# 4) Remove the two trailing paragraphs that followed it
#    ("我本不愿将青春往事..." and "最初的尝试...") entirely, including their
#    paragraph marks, leaving the final blank paragraph untouched.
$paraStart = $null
$paraEnd = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.StartsWith("我本不愿将青春往事翻出来")) {
        $paraStart = $p.Range.Start
    }
    if ($t.StartsWith("最初的尝试，不过出于本能")) {
        $paraEnd = $p.Range.End
    }
}

if (($paraStart -ne $null) -and ($paraEnd -ne $null)) {
    $killRange = $d.Range($paraStart, $paraEnd)
    $killRange.Delete()
}
